# Daily GitHub Actions refresh of the cryptos price list.
# Column D ("Price") cells hold text that often *looks* numeric
# (e.g. "548.91", "61.844.10", "0.0934"). Excel auto-coerces such
# strings to real numbers on a plain .Value assignment, which would
# silently rewrite the literal text (dropping trailing zeros, turning
# "7.00" into 7, mangling thousand-dot grouping like "61.844.10",
# etc). Forcing NumberFormat to Text ("@") before writing keeps these
# as the exact literal strings, matching the source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.844.10'
$ws.Range('E2').Value = '  -1.23%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.460.13'
$ws.Range('E3').Value = '  -2.93%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '548.91'
$ws.Range('E5').Value = '  -2.63%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.54'
$ws.Range('E6').Value = '  -1.76%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E8').Value = '  -3.88%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.457.09'
$ws.Range('E9').Value = '  -3.02%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.107'
$ws.Range('E10').Value = '  -4.44%  '
$ws.Range('E11').Value = '  -0.19%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.41'
$ws.Range('E12').Value = '  -1.86%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.353'
$ws.Range('E13').Value = '  -4.28%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.21'
$ws.Range('E14').Value = '  -2.93%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.902.57'
$ws.Range('E15').Value = '  -2.76%  '
$ws.Range('E16').Value = '  -2.01%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.647.42'
$ws.Range('E17').Value = '  -1.29%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.459.71'
$ws.Range('E18').Value = '  -2.21%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.96'
$ws.Range('E19').Value = '  -4.60%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.00'
$ws.Range('E20').Value = '  -3.18%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.17'
$ws.Range('E21').Value = '  -3.61%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '320.22'
$ws.Range('E22').Value = '  -2.67%  '
$ws.Range('E23').Value = '  +0.29%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.91'
$ws.Range('E24').Value = '  +4.74%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '63.95'
$ws.Range('E25').Value = '  -2.23%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0₃0988'
$ws.Range('E26').Value = '  -8.16%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.577.64'
$ws.Range('E27').Value = '  -2.61%  '
$ws.Range('E28').Value = '  +0.14%  '
$ws.Range('B29').Value = 'Fetch.AI'
$ws.Range('C29').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.48'
$ws.Range('E29').Value = '  -6.64%  '
$ws.Range('B30').Value = 'Bittensor'
$ws.Range('C30').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '534.36'
$ws.Range('E30').Value = '  -4.19%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.27'
$ws.Range('E31').Value = '  -5.62%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.75'
$ws.Range('E32').Value = '  -3.09%  '
$ws.Range('E33').Value = '  -5.38%  '
$ws.Range('E34').Value = '  -3.78%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.61'
$ws.Range('E35').Value = '  -1.92%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.76'
$ws.Range('E36').Value = '  -6.23%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  +0.06%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.79'
$ws.Range('E38').Value = '  -5.06%  '
$ws.Range('E39').Value = '  -2.43%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.30'
$ws.Range('E40').Value = '  -3.46%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.77'
$ws.Range('E41').Value = '  -0.38%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '140.65'
$ws.Range('E42').Value = '  -7.28%  '
$ws.Range('E43').Value = '  +0.08%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '40.37'
$ws.Range('E44').Value = '  -1.65%  '
$ws.Range('E45').Value = '  -4.61%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '143.60'
$ws.Range('E46').Value = '  -6.04%  '
$ws.Range('B47').Value = 'Filecoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.62'
$ws.Range('E47').Value = '  -3.23%  '
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '21.85'
$ws.Range('E48').Value = '  -3.08%  '
$ws.Range('E49').Value = '  -4.61%  '
$ws.Range('E50').Value = '  -2.63%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0934'
$ws.Range('E51').Value = '  -3.36%  '
